$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.909.18"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.25%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.742.35"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.34%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.34"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.19"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.14%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.601"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.53%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.14%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.381"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.52"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -17.73%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.228.06"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.43"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.564.23"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.55%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.747.84"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.03"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.73%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.80"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.84"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.74"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.536"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.63%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.93"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.36"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.82%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0884"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.86%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.90"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.71"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.18"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.07"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.77%  "

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.84"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.61%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.77"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.970"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.15"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +8.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.09"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "322.83"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.84"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.14"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0583"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.34%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.26"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.22%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.83"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.12%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.625"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.61%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.79%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.25%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.45%  "

